$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet had two header rows (row 1 and row 2) that together formed
# one logical header. Remove the second header row so all the data rows
# shift up by one, then rewrite row 1 as a single combined header row.
$ws.Rows("2").Delete()

# Rewrite the header row with the new column headers.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Match the font used elsewhere in the header/body (Arial 9) for the
# numeric-column headers.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

# Update the active selection to reflect the new header row, as in the
# target workbook.
$ws.Range("A2:K2").Select()
